# Generate Report for Handoff
#
# The localization status moved from "In Translation" to "Ready for
# handoff", and the associated "Latest Handoff"/"Latest HO Xliff Generate"
# timestamps were refreshed to reflect the handoff report regeneration.
# This touches the Overview roll-up sheet as well as the per-locale
# (zh-cn / de-de) detail sheets that back it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns + generation timestamp ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-31 02:43:57"

# --- zh-cn detail sheet: Status + Latest Handoff Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-31 02:43:52"

# --- de-de detail sheet: Status + Latest Handoff Datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-31 02:43:57"

# The new status text ("Ready for handoff") is longer than the old one
# ("In Translation"), so Excel's column autofit widens the Status columns
# on each touched sheet to keep the text from being clipped.
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null
